$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1282.174
$ws.Range("I28").Value = 869.5
$ws.Range("K28").Value = 869.5
$ws.Range("M28").Value = -384.5
$ws.Range("H98").Value = 4531.2104
$ws.Range("I98").Value = 4531.2104
$ws.Range("K98").Value = 4531.2104
$ws.Range("M98").Value = -3033.2104
$ws.Range("H107").Value = 759
$ws.Range("I107").Value = 689.6667
$ws.Range("K107").Value = 689.6667
$ws.Range("M107").Value = 1230.3333
$ws.Range("H116").Value = 3699.0908
$ws.Range("I116").Value = 3754.6667
$ws.Range("J116").Value = 3449
$ws.Range("K116").Value = 3754.6667
$ws.Range("L116").Value = 3449
$ws.Range("M116").Value = -312.6667000000002
$ws.Range("N116").Value = -10333
$ws.Range("H122").Value = 4531.2104
$ws.Range("I122").Value = 4531.2104
$ws.Range("K122").Value = 13593.6312
$ws.Range("M122").Value = -11143.6312
$ws.Range("H137").Value = 2322.9697
$ws.Range("I137").Value = 1904.8334
$ws.Range("J137").Value = 2824.7334
$ws.Range("K137").Value = 5714.5002
$ws.Range("L137").Value = 8474.200199999999
$ws.Range("M137").Value = -3164.5002
$ws.Range("N137").Value = -13574.2002
$ws.Range("H138").Value = 2915.28
$ws.Range("J138").Value = 3115.4246
$ws.Range("L138").Value = 9346.273799999999
$ws.Range("N138").Value = -19626.2738

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2741.81
$ws.Range("I32").Value = 2362.7217
$ws.Range("K32").Value = 2362.7217
$ws.Range("M32").Value = -2075.7217
$ws.Range("H45").Value = 16918.367
$ws.Range("J45").Value = 3323.2632
$ws.Range("L45").Value = 3323.2632
$ws.Range("N45").Value = -4077.2632
$ws.Range("H61").Value = 7775.2856
$ws.Range("I61").Value = 6798.6665
$ws.Range("J61").Value = 9533.200000000001
$ws.Range("K61").Value = 6798.6665
$ws.Range("L61").Value = 9533.200000000001
$ws.Range("M61").Value = -6586.6665
$ws.Range("N61").Value = -9957.200000000001
$ws.Range("H74").Value = 148911.66
$ws.Range("I74").Value = 200701.17
$ws.Range("J74").Value = 3901
$ws.Range("K74").Value = 200701.17
$ws.Range("L74").Value = 3901
$ws.Range("M74").Value = -199827.17
$ws.Range("N74").Value = -5649
$ws.Range("H77").Value = 148911.66
$ws.Range("I77").Value = 200701.17
$ws.Range("J77").Value = 3901
$ws.Range("K77").Value = 1003505.85
$ws.Range("L77").Value = 19505
$ws.Range("M77").Value = -999137.8500000001
$ws.Range("N77").Value = -28241
$ws.Range("H102").Value = 3442.0908
$ws.Range("I102").Value = 3252.641
$ws.Range("K102").Value = 3252.641
$ws.Range("M102").Value = -1630.641
$ws.Range("H136").Value = 7775.2856
$ws.Range("I136").Value = 6798.6665
$ws.Range("J136").Value = 9533.200000000001
$ws.Range("K136").Value = 20395.9995
$ws.Range("L136").Value = 28599.6
$ws.Range("M136").Value = -17845.9995
$ws.Range("N136").Value = -33699.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 53491.05
$ws.Range("I99").Value = 92759.17999999999
$ws.Range("K99").Value = 92759.17999999999
$ws.Range("M99").Value = -91261.17999999999
$ws.Range("H107").Value = 1545.5555
$ws.Range("I107").Value = 1545.5555
$ws.Range("K107").Value = 1545.5555
$ws.Range("M107").Value = 374.4445000000001
$ws.Range("H134").Value = 2477.3572
$ws.Range("I134").Value = 2161.125
$ws.Range("J134").Value = 4374.75
$ws.Range("K134").Value = 6483.375
$ws.Range("L134").Value = 13124.25
$ws.Range("M134").Value = -3948.375
$ws.Range("N134").Value = -18194.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66667690
$ws.Range("I7").Value = 1074.9166
$ws.Range("J7").Value = 333334140
$ws.Range("K7").Value = 1074.9166
$ws.Range("L7").Value = 333334140
$ws.Range("M7").Value = -961.9166
$ws.Range("N7").Value = -333334366
$ws.Range("H31").Value = 2830.47
$ws.Range("I31").Value = 2633.561
$ws.Range("J31").Value = 3727.5
$ws.Range("K31").Value = 2633.561
$ws.Range("L31").Value = 3727.5
$ws.Range("M31").Value = -2338.561
$ws.Range("N31").Value = -4317.5
$ws.Range("H34").Value = 2830.47
$ws.Range("I34").Value = 2633.561
$ws.Range("J34").Value = 3727.5
$ws.Range("K34").Value = 2633.561
$ws.Range("L34").Value = 3727.5
$ws.Range("M34").Value = -2431.561
$ws.Range("N34").Value = -4131.5
$ws.Range("H58").Value = 2968.6086
$ws.Range("J58").Value = 4191.8887
$ws.Range("L58").Value = 4191.8887
$ws.Range("N58").Value = -4597.8887
$ws.Range("H86").Value = 3107.7
$ws.Range("I86").Value = 2297.5715
$ws.Range("K86").Value = 2297.5715
$ws.Range("M86").Value = -1174.5715
$ws.Range("H89").Value = 3107.7
$ws.Range("I89").Value = 2297.5715
$ws.Range("K89").Value = 11487.8575
$ws.Range("M89").Value = -5871.8575
$ws.Range("H132").Value = 3839.7666
$ws.Range("I132").Value = 3630.9614
$ws.Range("J132").Value = 5197
$ws.Range("K132").Value = 10892.8842
$ws.Range("L132").Value = 15591
$ws.Range("M132").Value = -8362.8842
$ws.Range("N132").Value = -20651
$ws.Range("H136").Value = 2968.6086
$ws.Range("J136").Value = 4191.8887
$ws.Range("L136").Value = 12575.6661
$ws.Range("N136").Value = -17675.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2244.9546
$ws.Range("I129").Value = 2109.818
$ws.Range("J129").Value = 2380.0908
$ws.Range("K129").Value = 6329.454000000001
$ws.Range("L129").Value = 7140.2724
$ws.Range("M129").Value = -1329.454000000001
$ws.Range("N129").Value = -17140.2724
$ws.Range("H139").Value = 3841.842
$ws.Range("J139").Value = 3999.6924
$ws.Range("L139").Value = 11999.0772
$ws.Range("N139").Value = -22279.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5539.2256
$ws.Range("I113").Value = 5791.478
$ws.Range("J113").Value = 4814
$ws.Range("K113").Value = 5791.478
$ws.Range("L113").Value = 4814
$ws.Range("M113").Value = -3621.478
$ws.Range("N113").Value = -9154
$ws.Range("H132").Value = 5080.3
$ws.Range("I132").Value = 4534.3335
$ws.Range("J132").Value = 9994
$ws.Range("K132").Value = 13603.0005
$ws.Range("L132").Value = 29982
$ws.Range("M132").Value = -11073.0005
$ws.Range("N132").Value = -35042

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7883.9375
$ws.Range("I7").Value = 6321
$ws.Range("K7").Value = 6321
$ws.Range("M7").Value = -6209
$ws.Range("H40").Value = 21308.705
$ws.Range("I40").Value = 35676.656
$ws.Range("K40").Value = 35676.656
$ws.Range("M40").Value = -35540.656
$ws.Range("H46").Value = 1461.659
$ws.Range("I46").Value = 1835.8572
$ws.Range("K46").Value = 1835.8572
$ws.Range("M46").Value = -1647.8572
$ws.Range("H55").Value = 632.6923
$ws.Range("I55").Value = 561.5714
$ws.Range("J55").Value = 715.6667
$ws.Range("K55").Value = 561.5714
$ws.Range("L55").Value = 715.6667
$ws.Range("M55").Value = -388.5714
$ws.Range("N55").Value = -1061.6667
$ws.Range("H68").Value = 4380.4
$ws.Range("I68").Value = 3634
$ws.Range("J68").Value = 5500
$ws.Range("K68").Value = 3634
$ws.Range("L68").Value = 5500
$ws.Range("M68").Value = -2885
$ws.Range("N68").Value = -6998
$ws.Range("H71").Value = 4380.4
$ws.Range("I71").Value = 3634
$ws.Range("J71").Value = 5500
$ws.Range("K71").Value = 18170
$ws.Range("L71").Value = 27500
$ws.Range("M71").Value = -14426
$ws.Range("N71").Value = -34988
$ws.Range("H82").Value = 5805.4287
$ws.Range("I82").Value = 4971.778
$ws.Range("K82").Value = 4971.778
$ws.Range("M82").Value = -4610.778
$ws.Range("H85").Value = 5805.4287
$ws.Range("I85").Value = 4971.778
$ws.Range("K85").Value = 4971.778
$ws.Range("M85").Value = -3723.778
$ws.Range("H126").Value = 7883.9375
$ws.Range("I126").Value = 6321
$ws.Range("K126").Value = 18963
$ws.Range("M126").Value = -16493
$ws.Range("H132").Value = 4599.9
$ws.Range("I132").Value = 4199.8
$ws.Range("K132").Value = 12599.4
$ws.Range("M132").Value = -10069.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 554.5
$ws.Range("I107").Value = 555.4375
$ws.Range("K107").Value = 1666.3125
$ws.Range("M107").Value = 253.6875
$ws.Range("H113").Value = 982.5789
$ws.Range("I113").Value = 1179.5555
$ws.Range("K113").Value = 3538.6665
$ws.Range("M113").Value = -1368.6665
$ws.Range("H126").Value = 2097.9375
$ws.Range("I126").Value = 2097.9375
$ws.Range("K126").Value = 6293.8125
$ws.Range("M126").Value = -3823.8125
$ws.Range("H132").Value = 4055.818
$ws.Range("I132").Value = 3751.4792
$ws.Range("J132").Value = 6142.7144
$ws.Range("K132").Value = 11254.4376
$ws.Range("L132").Value = 18428.1432
$ws.Range("M132").Value = -8724.437600000001
$ws.Range("N132").Value = -23488.1432

$wb.Save()